$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction on Login tests: fix forbidden-character test data
$ws.Range("B4").Value = "/""&)ç_'"
$ws.Range("C4").Value = "/""&)ç_'"
$ws.Range("A5").Value = "([ç^"

# Update the active selection to match the author's saved view state
$ws.Range("D17").Select()
